$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Course_packages: bump F4, append 3 new package rows (6-8), new zoom/selection
# ---------------------------------------------------------------------------
$wsPkg = $wb.Worksheets.Item("Course_packages")
$wsPkg.Activate()

# Existing row 4 correction (F4 233 -> 308)
$wsPkg.Range("F4").Value = 308

# Append new rows 6,7,8 - copy formatting from row 5 (the last existing data row)
# so the new cells pick up the same number/date styles instead of minting new ones.
$wsPkg.Range("A5:F5").Copy()
$wsPkg.Range("A6:F8").PasteSpecial(-4122)

$wsPkg.Range("A6").Value = 5
$wsPkg.Range("B6").Value = 10
$wsPkg.Range("C6").Value = 44228
$wsPkg.Range("D6").Value = 44531
$wsPkg.Range("E6").Value = "HAPPY NEW YEAR"
$wsPkg.Range("F6").Value = 188

$wsPkg.Range("A7").Value = 6
$wsPkg.Range("B7").Value = 10
$wsPkg.Range("C7").Value = 44256
$wsPkg.Range("D7").Value = 44532
$wsPkg.Range("E7").Value = "HAPPY NEW MONTH"
$wsPkg.Range("F7").Value = 198

$wsPkg.Range("A8").Value = 7
$wsPkg.Range("B8").Value = 15
$wsPkg.Range("C8").Value = 44287
$wsPkg.Range("D8").Value = 44533
$wsPkg.Range("E8").Value = "HAPPY APRIL FOOLS"
$wsPkg.Range("F8").Value = 233

$excel.ActiveWindow.Zoom = 150
$wsPkg.Range("G10").Select()

# ---------------------------------------------------------------------------
# 2) Buys: clean up remaining_redemptions, append 4 new rows
# ---------------------------------------------------------------------------
$wsBuys = $wb.Worksheets.Item("Buys")
$wsBuys.Activate()

$wsBuys.Range("D2").Value = 4
$wsBuys.Range("D3").Value = 3
$wsBuys.Range("D6").Value = 18

# Append new rows 8-11 - copy formatting from row 7 (the last existing data row)
$wsBuys.Range("A7:D7").Copy()
$wsBuys.Range("A8:D11").PasteSpecial(-4122)

$wsBuys.Range("A8").Value = 5
$wsBuys.Range("B8").Value = "A123456789014"
$wsBuys.Range("C8").Value = 44288
$wsBuys.Range("D8").Value = 10

$wsBuys.Range("A9").Value = 5
$wsBuys.Range("B9").Value = "A123456789020"
$wsBuys.Range("C9").Value = 44288
$wsBuys.Range("D9").Value = 10

$wsBuys.Range("A10").Value = 6
$wsBuys.Range("B10").Value = "A123456789033"
$wsBuys.Range("C10").Value = 44288
$wsBuys.Range("D10").Value = 10

$wsBuys.Range("A11").Value = 7
$wsBuys.Range("B11").Value = "A123456789019"
$wsBuys.Range("C11").Value = 44288
$wsBuys.Range("D11").Value = 15

# ---------------------------------------------------------------------------
# 3) Pure view/selection updates on other sheets touched while testing
# ---------------------------------------------------------------------------
$wsSessions = $wb.Worksheets.Item("Sessions")
$wsSessions.Activate()
$wsSessions.Range("D10").Select()

$wsCredit = $wb.Worksheets.Item("Credit_cards")
$wsCredit.Activate()
$wsCredit.Range("A6").Select()

$wsOwns = $wb.Worksheets.Item("Owns")
$wsOwns.Activate()
$wsOwns.Range("B9").Select()

$wsRedeems = $wb.Worksheets.Item("Redeems")
$wsRedeems.Activate()
$wsRedeems.Range("B3").Select()

$wsCancels = $wb.Worksheets.Item("Cancels")
$wsCancels.Activate()
$wsCancels.Range("H15").Select()

# Buys ends up as the active/selected sheet (matches the final activeTab in the diff),
# so select it last.
$wsBuys.Activate()
$wsBuys.Range("D13").Select()
